$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 226
$ws1.Range("F4").Value = 4794
$ws1.Range("F5").Value = 0
$ws1.Range("F7").Value = 113
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 0
$ws1.Range("F10").Value = 746
$ws1.Range("F11").Value = 220
$ws1.Range("F12").Value = 1164
$ws1.Range("F14").Value = 259
$ws1.Range("F15").Value = 181
$ws1.Range("F16").Value = 82
$ws1.Range("F20").Value = 0
$ws1.Range("F21").Value = 0
$ws1.Range("F23").Value = 0
$ws1.Range("F24").Value = 539
$ws1.Range("F25").Value = 48
$ws1.Range("F27").Value = 0
$ws1.Range("F28").Value = 39
$ws1.Range("F29").Value = 12
$ws1.Range("F30").Value = 2566
$ws1.Range("F31").Value = 568
$ws1.Range("F32").Value = 0
$ws1.Range("F36").Value = 0
$ws1.Range("F37").Value = 170
$ws1.Range("F39").Value = 1559
$ws1.Range("F40").Value = 0
$ws1.Range("F42").Value = 63
$ws1.Range("F43").Value = 57
$ws1.Range("F47").Value = 0

# Sheet "全部类型" (All types) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 110
$ws4.Range("F5").Value = 210
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 113
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 93
$ws4.Range("F11").Value = 746
$ws4.Range("F12").Value = 220
$ws4.Range("F16").Value = 181
$ws4.Range("F17").Value = 0
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 6279
$ws4.Range("F25").Value = 539
$ws4.Range("F27").Value = 3966
$ws4.Range("F29").Value = 39
$ws4.Range("F30").Value = 12
$ws4.Range("F32").Value = 568
$ws4.Range("F33").Value = 528
$ws4.Range("F34").Value = 139
$ws4.Range("F35").Value = 0
$ws4.Range("F37").Value = 369
$ws4.Range("F38").Value = 170
$ws4.Range("F39").Value = 7
$ws4.Range("F40").Value = 0
$ws4.Range("F41").Value = 959
$ws4.Range("F43").Value = 63
$ws4.Range("F44").Value = 57
$ws4.Range("F45").Value = 493
$ws4.Range("F47").Value = 0
$ws4.Range("F48").Value = 582
